$d = $word.ActiveDocument

# 1) Append " +/-" as a new run at the end of the "equivalencias" paragraph.
#    Toggling a character property (and reverting it) forces the engine to
#    split off a fresh run instead of merging the inserted text into the
#    preceding run, while leaving the run's rPr identical to its neighbour.
$p2 = $d.Paragraphs.Item(2)
$endRange = $p2.Range
$endRange.SetRange($endRange.End - 1, $endRange.End - 1)
$endRange.InsertAfter(" +/-")
$endRange.Font.Bold = 1
$endRange.Font.Bold = 0

# 2) Clear the "Preguntar por el practico de Programación I" paragraph,
#    leaving an empty paragraph behind.
$found = $d.Content.Find.Execute("Preguntar por el practico de Programación I", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# 3) Clear the "Preguntar por TP N°2 (Arquitectura de las computadoras) en horario de consulta." paragraph,
#    leaving an empty paragraph behind.
$found = $d.Content.Find.Execute("Preguntar por TP N°2 (Arquitectura de las computadoras) en horario de consulta.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
